$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.127881588408715
$ws.Range("C2").Value = 1.667794583268128
$ws.Range("D2").Value = 26.21740644021617
$ws.Range("E2").Value = 616238.5361209477
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 616266.5492035595
